$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 482.14285
$ws.Cells.Item(41, 9).Value = 459.18182
$ws.Cells.Item(41, 10).Value = 566.3333
$ws.Cells.Item(41, 11).Value = 459.18182
$ws.Cells.Item(41, 12).Value = 566.3333
$ws.Cells.Item(41, 13).Value = -19.18182000000002
$ws.Cells.Item(41, 14).Value = -1446.3333
$ws.Cells.Item(51, 8).Value = 3386.1724
$ws.Cells.Item(51, 9).Value = 4334
$ws.Cells.Item(51, 10).Value = 3234.52
$ws.Cells.Item(51, 11).Value = 4334
$ws.Cells.Item(51, 12).Value = 3234.52
$ws.Cells.Item(51, 13).Value = -3850
$ws.Cells.Item(51, 14).Value = -4202.52
$ws.Cells.Item(80, 8).Value = 2995.4443
$ws.Cells.Item(80, 9).Value = 2472
$ws.Cells.Item(80, 10).Value = 3414.2
$ws.Cells.Item(80, 11).Value = 7416
$ws.Cells.Item(80, 12).Value = 10242.6
$ws.Cells.Item(80, 13).Value = -6418
$ws.Cells.Item(80, 14).Value = -12238.6
$ws.Cells.Item(83, 8).Value = 2995.4443
$ws.Cells.Item(83, 9).Value = 2472
$ws.Cells.Item(83, 10).Value = 3414.2
$ws.Cells.Item(83, 11).Value = 22248
$ws.Cells.Item(83, 12).Value = 30727.8
$ws.Cells.Item(83, 13).Value = -17256
$ws.Cells.Item(83, 14).Value = -40711.8
$ws.Cells.Item(132, 8).Value = 4097.9546
$ws.Cells.Item(132, 9).Value = 3656.2942
$ws.Cells.Item(132, 10).Value = 5599.6
$ws.Cells.Item(132, 11).Value = 10968.8826
$ws.Cells.Item(132, 12).Value = 16798.8
$ws.Cells.Item(132, 13).Value = -8438.882599999999
$ws.Cells.Item(132, 14).Value = -21858.8
$ws.Cells.Item(138, 8).Value = 5704.8887
$ws.Cells.Item(138, 9).Value = 9299.25
$ws.Cells.Item(138, 10).Value = 4677.9287
$ws.Cells.Item(138, 11).Value = 27897.75
$ws.Cells.Item(138, 12).Value = 14033.7861
$ws.Cells.Item(138, 13).Value = -22757.75
$ws.Cells.Item(138, 14).Value = -24313.7861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1570.1404
$ws.Cells.Item(32, 9).Value = 1030.0834
$ws.Cells.Item(32, 10).Value = 4450.4443
$ws.Cells.Item(32, 11).Value = 1030.0834
$ws.Cells.Item(32, 12).Value = 4450.4443
$ws.Cells.Item(32, 13).Value = -743.0834
$ws.Cells.Item(32, 14).Value = -5024.4443
$ws.Cells.Item(45, 8).Value = 2293.1738
$ws.Cells.Item(45, 9).Value = 928.4375
$ws.Cells.Item(45, 10).Value = 5412.5713
$ws.Cells.Item(45, 11).Value = 928.4375
$ws.Cells.Item(45, 12).Value = 5412.5713
$ws.Cells.Item(45, 13).Value = -551.4375
$ws.Cells.Item(45, 14).Value = -6166.5713
$ws.Cells.Item(74, 8).Value = 2767.2273
$ws.Cells.Item(74, 9).Value = 1669.1538
$ws.Cells.Item(74, 10).Value = 4353.3335
$ws.Cells.Item(74, 11).Value = 1669.1538
$ws.Cells.Item(74, 12).Value = 4353.3335
$ws.Cells.Item(74, 13).Value = -795.1538
$ws.Cells.Item(74, 14).Value = -6101.3335
$ws.Cells.Item(77, 8).Value = 2767.2273
$ws.Cells.Item(77, 9).Value = 1669.1538
$ws.Cells.Item(77, 10).Value = 4353.3335
$ws.Cells.Item(77, 11).Value = 8345.769
$ws.Cells.Item(77, 12).Value = 21766.6675
$ws.Cells.Item(77, 13).Value = -3977.769
$ws.Cells.Item(77, 14).Value = -30502.6675
$ws.Cells.Item(102, 8).Value = 1373.7273
$ws.Cells.Item(102, 9).Value = 1373.7273
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1373.7273
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = 248.2727
$ws.Cells.Item(132, 8).Value = 3039.6365
$ws.Cells.Item(132, 9).Value = 2473.7
$ws.Cells.Item(132, 10).Value = 8699
$ws.Cells.Item(132, 11).Value = 7421.099999999999
$ws.Cells.Item(132, 12).Value = 26097
$ws.Cells.Item(132, 13).Value = -4891.099999999999
$ws.Cells.Item(132, 14).Value = -31157

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1104.75
$ws.Cells.Item(105, 9).Value = 1134.0714
$ws.Cells.Item(105, 10).Value = 899.5
$ws.Cells.Item(105, 11).Value = 1134.0714
$ws.Cells.Item(105, 12).Value = 899.5
$ws.Cells.Item(105, 13).Value = 612.9286
$ws.Cells.Item(105, 14).Value = -4393.5
$ws.Cells.Item(134, 8).Value = 5863.8237
$ws.Cells.Item(134, 9).Value = 4715.2593
$ws.Cells.Item(134, 10).Value = 10294
$ws.Cells.Item(134, 11).Value = 14145.7779
$ws.Cells.Item(134, 12).Value = 30882
$ws.Cells.Item(134, 13).Value = -11610.7779
$ws.Cells.Item(134, 14).Value = -35952

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5296.115
$ws.Cells.Item(31, 9).Value = 4280.2
$ws.Cells.Item(31, 10).Value = 5931.0625
$ws.Cells.Item(31, 11).Value = 4280.2
$ws.Cells.Item(31, 12).Value = 5931.0625
$ws.Cells.Item(31, 13).Value = -3985.2
$ws.Cells.Item(31, 14).Value = -6521.0625
$ws.Cells.Item(34, 8).Value = 5296.115
$ws.Cells.Item(34, 9).Value = 4280.2
$ws.Cells.Item(34, 10).Value = 5931.0625
$ws.Cells.Item(34, 11).Value = 4280.2
$ws.Cells.Item(34, 12).Value = 5931.0625
$ws.Cells.Item(34, 13).Value = -4078.2
$ws.Cells.Item(34, 14).Value = -6335.0625
$ws.Cells.Item(134, 8).Value = 4247.0625
$ws.Cells.Item(134, 9).Value = 4125.25
$ws.Cells.Item(134, 10).Value = 5099.75
$ws.Cells.Item(134, 11).Value = 12375.75
$ws.Cells.Item(134, 12).Value = 15299.25
$ws.Cells.Item(134, 13).Value = -9840.75
$ws.Cells.Item(134, 14).Value = -20369.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 4409.778
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 4409.778
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 13229.334
$ws.Cells.Item(68, 14).Value = -14851.334
$ws.Cells.Item(71, 8).Value = 4409.778
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 4409.778
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 39688.002
$ws.Cells.Item(71, 14).Value = -47800.002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1927.5834
$ws.Cells.Item(122, 9).Value = 1863.2
$ws.Cells.Item(122, 10).Value = 2249.5
$ws.Cells.Item(122, 11).Value = 5589.6
$ws.Cells.Item(122, 12).Value = 6748.5
$ws.Cells.Item(122, 13).Value = -3139.6
$ws.Cells.Item(122, 14).Value = -11648.5
$ws.Cells.Item(132, 8).Value = 3732.1904
$ws.Cells.Item(132, 9).Value = 3868.9
$ws.Cells.Item(132, 10).Value = 998
$ws.Cells.Item(132, 11).Value = 11606.7
$ws.Cells.Item(132, 12).Value = 2994
$ws.Cells.Item(132, 13).Value = -9076.700000000001
$ws.Cells.Item(132, 14).Value = -8054

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1865.5
$ws.Cells.Item(22, 9).Value = 1162.25
$ws.Cells.Item(22, 10).Value = 2568.75
$ws.Cells.Item(22, 11).Value = 1162.25
$ws.Cells.Item(22, 12).Value = 2568.75
$ws.Cells.Item(22, 13).Value = -867.25
$ws.Cells.Item(22, 14).Value = -3158.75
$ws.Cells.Item(27, 8).Value = 1865.5
$ws.Cells.Item(27, 9).Value = 1162.25
$ws.Cells.Item(27, 10).Value = 2568.75
$ws.Cells.Item(27, 11).Value = 1162.25
$ws.Cells.Item(27, 12).Value = 2568.75
$ws.Cells.Item(27, 13).Value = -1055.25
$ws.Cells.Item(27, 14).Value = -2782.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2950.25
$ws.Cells.Item(126, 9).Value = 2933.6667
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 8801.000100000001
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -6331.000100000001
$ws.Cells.Item(126, 14).Value = -13940
$ws.Cells.Item(132, 8).Value = 2168
$ws.Cells.Item(132, 9).Value = 1357.5714
$ws.Cells.Item(132, 10).Value = 2640.75
$ws.Cells.Item(132, 11).Value = 4072.7142
$ws.Cells.Item(132, 12).Value = 7922.25
$ws.Cells.Item(132, 13).Value = -1542.7142
$ws.Cells.Item(132, 14).Value = -12982.25
$ws.Cells.Item(136, 8).Value = 6200.3887
$ws.Cells.Item(136, 9).Value = 5720.857
$ws.Cells.Item(136, 10).Value = 10899.8
$ws.Cells.Item(136, 11).Value = 17162.571
$ws.Cells.Item(136, 12).Value = 32699.4
$ws.Cells.Item(136, 13).Value = -14612.571
$ws.Cells.Item(136, 14).Value = -37799.39999999999
